# Add a "Save" column (column H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (G1) onto the new
# header cell (H1) so it matches the other bold/bordered header cells,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" values for each data row (unstyled, numeric 1).
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
